$wb = $excel.ActiveWorkbook

# --- Locate existing sheets ---
$wsEstandar = $wb.Worksheets.Item("Estandar (ilerdair)")

# --- Insert the new "RGBN Pruebas" sheet right after "Estandar (ilerdair)" ---
$wsRgbn = $wb.Worksheets.Add($null, $wsEstandar)
$wsRgbn.Name = "RGBN Pruebas"

# --- Fill in the new sheet's data (mirrors the layout of the other config sheets) ---
$wsRgbn.Range("A1").Value = "master"
$wsRgbn.Range("A2").Value = 1
$wsRgbn.Range("A3").Value = 2
$wsRgbn.Range("A4").Value = 3
$wsRgbn.Range("A5").Value = 4
$wsRgbn.Range("A6").Value = 5

$wsRgbn.Range("B1").Value = "R780"
$wsRgbn.Range("B2").Value = "R450"
$wsRgbn.Range("B3").Value = "R550"
$wsRgbn.Range("B4").Value = "R670"
$wsRgbn.Range("B5").Value = "R710"
$wsRgbn.Range("B6").Value = "R800"

# --- Reset the selection on the other two config sheets to the full data range ---
# (re-fetch by name: inserting the new sheet shifts everyone's position, and
#  cached worksheet references can otherwise resolve against a stale index)
$wsConfig1 = $wb.Worksheets.Item("Config 1 Algerri")
$wsConfig1.Activate() | Out-Null
$wsConfig1.Range("A1:B6").Select() | Out-Null

$wsConfig2 = $wb.Worksheets.Item("Config 2 Algerri")
$wsConfig2.Activate() | Out-Null
$wsConfig2.Range("A1:B6").Select() | Out-Null

# --- Make the new RGBN Pruebas sheet the active one, with B5 selected ---
$wsRgbn.Activate() | Out-Null
$wsRgbn.Range("B5").Select() | Out-Null
